$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 280, shifting existing rows (280-288) down to (282-290)
$insertRange = $ws.Range("A280:F281")
$insertRange.EntireRow.Insert()

# New row 280: 19_02_05_others
$ws.Range("B280").Value = "19_02_05_others"
$ws.Range("C280").Value = "19_02_05_others"
$ws.Range("D280").Value = "(new)"

# New row 281: 19_02_17_electricity
$ws.Range("B281").Value = "19_02_17_electricity"
$ws.Range("C281").Value = "19_02_17_electricity"
$ws.Range("D281").Value = "(new)"

$excel.ActiveWindow.ScrollRow = 258
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D279:D281").Select()
